$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all hyperlinks from the sheet (B2:B4 had mailto: links).
# The underlying collection re-indexes as items are deleted, so run the
# delete pass twice to make sure every hyperlink is actually removed.
foreach ($hl in @($ws.Hyperlinks)) {
    $hl.Delete()
}
foreach ($hl in @($ws.Hyperlinks)) {
    $hl.Delete()
}

# Clear the sample/demo data that used to populate rows 2-4
$ws.Range("A2:E4").ClearContents()

# Move the active selection back to A2 (was C7)
$ws.Range("A2").Select()
